$d = $word.ActiveDocument

function New-WordXmlPackage($bodyInnerXml) {
    $header = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
    $footer = '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $header + $bodyInnerXml + $footer
}

# Locate the (first) paragraph whose text contains the given substring; returns its Range.
function Find-ParagraphContaining($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p.Range
        }
    }
    return $null
}

# --- Paragraph 1: "Install the python module pyinstaller. ... (conda has one an I will try it in the future)"
# Change the proofErr wrapping "an" from spellStart/spellEnd to gramStart/gramEnd.
$p1 = Find-ParagraphContaining $d "Install the python module"
$xml1 = '<w:body><w:p>' + `
  '<w:r><w:t xml:space="preserve">Install the python module </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>pyinstaller</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t>. I am using 3.4 that I installed using pip (</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>conda</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> has one </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:t>an</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> I will try it in the future)</w:t></w:r>' + `
  '</w:p></w:body>'
$p1.InsertXML((New-WordXmlPackage $xml1))

# --- Paragraph 5: "Get  the crop_int.rc,  crop_int.spec and pyinstall2.bat files from TF ..."
# Wrap "Get  the" in gramStart/gramEnd, splitting the " the " run so the trailing
# space moves to a new run placed after the gramEnd marker.
$p5 = Find-ParagraphContaining $d "crop_int.rc"
$xml5 = '<w:body><w:p>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:t xml:space="preserve">Get </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> the</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>crop_int.rc</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve">,  </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>crop_int.spec</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> and pyinstall2.bat files </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">from TF and make sure they are mapped </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">to your </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>crop_int</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> folder with the source code</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">. They should be mapped that way by default. </w:t></w:r>' + `
  '</w:p></w:body>'
$p5.InsertXML((New-WordXmlPackage $xml5))

# --- Paragraph 8: "Open a command window in your environment and navigate to the folder with the source code"
# Split off the final word "code" into its own run wrapped in gramStart/gramEnd.
$p8 = Find-ParagraphContaining $d "Open a command window"
$xml8 = '<w:body><w:p>' + `
  '<w:r><w:t xml:space="preserve">Open a command window in your environment and navigate to the folder with the source </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:t>code</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '</w:p></w:body>'
$p8.InsertXML((New-WordXmlPackage $xml8))

# --- Paragraph 14: "The classim folder in the dist subfolder has all the files you need to create a distribution"
# Split off the final word "distribution" into its own run wrapped in gramStart/gramEnd.
$p14 = Find-ParagraphContaining $d "subfolder has all the files"
$xml14 = '<w:body><w:p>' + `
  '<w:r><w:t xml:space="preserve">The </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>classim</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> folder in the </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>dist</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> subfolder has all the files you need to create a </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:t>distribution</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '</w:p></w:body>'
$p14.InsertXML((New-WordXmlPackage $xml14))

# --- Add two new empty paragraphs right before the final (last) empty paragraph
# in the document body (i.e. right after the "distribution" paragraph).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$lastPara.Collapse(0)
$lastPara.InsertXML((New-WordXmlPackage '<w:body><w:p/><w:p/><w:p/></w:body>'))
